$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values are plain numeric
# strings (e.g. "0.9994", "86.23") so Excel/COM does not silently coerce
# them into numeric cells and lose the original text representation.
$ws.Range("D2").Value = "29.969.31"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.877.01"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.57"
$ws.Range("E5").Value = "  -3.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4928"
$ws.Range("E7").Value = "  -3.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2922"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06638"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").Value = "1.877.07"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.73"
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07241"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6664"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.23"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.878"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "29.953.95"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007873"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "2.120.73"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.777"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.783"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.057"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.76"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.06"
$ws.Range("E26").Value = "  +5.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.04"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.915"
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.390"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.194"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08731"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.968"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05062"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7121"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.116"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01784"
$ws.Range("E37").Value = "  +5.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.687"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.178"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9317"
$ws.Range("E40").Value = "  -3.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4249"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9984"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.756"
$ws.Range("E43").Value = "  -6.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.75"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.422"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1270"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05660"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.67"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.314"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3786"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.95"
